# Weekly update: insert the new week's Melón price data (rows 30-32) and
# shift the existing history down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 30 (pushes former rows 30-51 down to 33-54).
$ws.Rows("30:32").Insert()

# New row 30 - Melón, Tuna, Extra - week of 2021-12-24
$ws.Cells.Item(30, 1).Value = 8
$ws.Cells.Item(30, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44554
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = 100112027
$ws.Cells.Item(30, 7).Value = "Melón"
$ws.Cells.Item(30, 8).Value = "Tuna"
$ws.Cells.Item(30, 9).Value = "Extra"
$ws.Cells.Item(30, 10).Value = 4000
$ws.Cells.Item(30, 11).Value = 1000
$ws.Cells.Item(30, 12).Value = 1100
$ws.Cells.Item(30, 13).Value = 1050
$ws.Cells.Item(30, 14).Value = "$/unidad"
$ws.Cells.Item(30, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(30, 16).Value = 1050
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = "Hortaliza"

# New row 31 - Melón, Tuna, Primera - week of 2021-12-24
$ws.Cells.Item(31, 1).Value = 8
$ws.Cells.Item(31, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 44554
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112027
$ws.Cells.Item(31, 7).Value = "Melón"
$ws.Cells.Item(31, 8).Value = "Tuna"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 4000
$ws.Cells.Item(31, 11).Value = 800
$ws.Cells.Item(31, 12).Value = 900
$ws.Cells.Item(31, 13).Value = 850
$ws.Cells.Item(31, 14).Value = "$/unidad"
$ws.Cells.Item(31, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(31, 16).Value = 850
$ws.Cells.Item(31, 17).Value = 1
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# New row 32 - Melón, Tuna, Super - week of 2021-12-24
$ws.Cells.Item(32, 1).Value = 8
$ws.Cells.Item(32, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44554
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = 100112027
$ws.Cells.Item(32, 7).Value = "Melón"
$ws.Cells.Item(32, 8).Value = "Tuna"
$ws.Cells.Item(32, 9).Value = "Super"
$ws.Cells.Item(32, 10).Value = 5000
$ws.Cells.Item(32, 11).Value = 1300
$ws.Cells.Item(32, 12).Value = 1400
$ws.Cells.Item(32, 13).Value = 1350
$ws.Cells.Item(32, 14).Value = "$/unidad"
$ws.Cells.Item(32, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(32, 16).Value = 1350
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = "Hortaliza"
